# Continue filling in the "Clientes" / "Dados Manjerico" workbook:
# split the single "Endereco" column into Numero / Complemento / Bairro,
# and append CEP / Telefone columns on both sheets, with the new data
# that was gathered for each customer.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Clientes
$ws2 = $wb.Worksheets.Item(2)   # Dados Manjerico

# ---------------------------------------------------------------
# Sheet 1: "Clientes"
#   before: A Identificador | B Razao Social | C CNPJ | D Endereco | E Cidade | F Estado
#   after : A Identificador | B Razao Social | C CNPJ | D Endereco | E Numero | F Complemento | G Bairro | H Cidade | I Estado | J CEP | K Telefone
# ---------------------------------------------------------------

# Insert three fresh columns (Numero, Complemento, Bairro) right before the
# existing "Cidade" column (E), pushing Cidade/Estado to H/I.
$ws1.Range("E1:G1").EntireColumn.Insert()

$ws1.Cells.Item(1,5).Value = "Número"
$ws1.Cells.Item(1,6).Value = "Complemento"
$ws1.Cells.Item(1,7).Value = "Bairro"

# Two more columns tacked on the end for CEP / Telefone.
$ws1.Cells.Item(1,10).Value = "CEP"
$ws1.Cells.Item(1,11).Value = "Telefone"

# Row 2 - Eventos & Cia
$ws1.Cells.Item(2,4).Value = "Rua Presidente Kennedy"
$ws1.Cells.Item(2,5).Value = 57
$ws1.Cells.Item(2,6).Value = "Casa 05"
$ws1.Cells.Item(2,7).Value = "Cônego"
$ws1.Cells.Item(2,10).Value = 28621000
$ws1.Cells.Item(2,11).Value = "(22) 2522-5120"

# Row 3 - Petrobras
$ws1.Cells.Item(3,4).Value = "Avenida República do Chile"
$ws1.Cells.Item(3,5).Value = 65
$ws1.Cells.Item(3,7).Value = "Centro"
$ws1.Cells.Item(3,10).Value = 20031912
$ws1.Cells.Item(3,11).Value = "(22) 2513-0056"

# Column widths (approximate best-fit sizing for the new columns)
$ws1.Columns.Item(2).ColumnWidth = 27.498697916666668
$ws1.Columns.Item(4).ColumnWidth = 22.498697916666668
$ws1.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws1.Columns.Item(6).ColumnWidth = 12.166666666666666
$ws1.Columns.Item(7).ColumnWidth = 6.498697916666667
$ws1.Columns.Item(8).ColumnWidth = 11.998697916666666
$ws1.Columns.Item(9).ColumnWidth = 5.998697916666667
$ws1.Columns.Item(10).ColumnWidth = 8.330729166666666
$ws1.Columns.Item(11).ColumnWidth = 12.830729166666666

# Make the identifier header bold & black (new font style picked up while typing).
$ws1.Range("A1").Font.Bold = $true
$ws1.Range("A1").Font.Color = 0

$ws1.Range("G4").Select() | Out-Null

# ---------------------------------------------------------------
# Sheet 2: "Dados Manjerico"
#   before: A Razao Social | B CNPJ | C Endereco | D Cidade | E Estado
#   after : A Razao Social | B CNPJ | C Endereco | D Numero | E Complemento | F Bairro | G Cidade | H Estado | I CEP | J Telefone
# ---------------------------------------------------------------

$ws2.Range("D1:F1").EntireColumn.Insert()

$ws2.Cells.Item(1,4).Value = "Número"
$ws2.Cells.Item(1,5).Value = "Complemento"
$ws2.Cells.Item(1,6).Value = "Bairro"

$ws2.Cells.Item(1,9).Value = "CEP"
$ws2.Cells.Item(1,10).Value = "Telefone"

$ws2.Cells.Item(2,4).Value = 26
$ws2.Cells.Item(2,5).Value = "Apto 212"
$ws2.Cells.Item(2,6).Value = "Tijuca"
$ws2.Cells.Item(2,9).Value = 20510150
$ws2.Cells.Item(2,10).Value = "(21) 2135-1448"

$ws2.Columns.Item(4).ColumnWidth = 23.166666666666668
$ws2.Columns.Item(5).ColumnWidth = 23.166666666666668
$ws2.Columns.Item(6).ColumnWidth = 23.166666666666668
$ws2.Columns.Item(10).ColumnWidth = 12.830729166666666

$ws2.Range("B39").Select() | Out-Null
$ws2.Activate() | Out-Null
